$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.309.55'
$ws.Range('E2').Value = '  -0.46%  '
$ws.Range('D3').Value = '2.462.32'
$ws.Range('E3').Value = '  +0.91%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '577.12'
$ws.Range('E5').Value = '  +1.19%  '
$ws.Range('D6').Value = '144.50'
$ws.Range('E6').Value = '  -0.40%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '0.533'
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').Value = '2.461.18'
$ws.Range('E9').Value = '  +1.04%  '
$ws.Range('D10').Value = '0.109'
$ws.Range('E10').Value = '  -0.57%  '
$ws.Range('E11').Value = '  +2.05%  '
$ws.Range('D12').Value = '5.23'
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('D13').Value = '0.346'
$ws.Range('E13').Value = '  -2.61%  '
$ws.Range('D14').Value = '26.46'
$ws.Range('E14').Value = '  -2.05%  '
$ws.Range('D15').Value = '0.0000175'
$ws.Range('E15').Value = '  -0.84%  '
$ws.Range('D16').Value = '2.902.08'
$ws.Range('E16').Value = '  +0.60%  '
$ws.Range('D17').Value = '62.101.62'
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('D18').Value = '2.457.57'
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('D19').Value = '10.89'
$ws.Range('E19').Value = '  -3.33%  '
$ws.Range('E20').Value = '  -1.00%  '
$ws.Range('D21').Value = '328.62'
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('D22').Value = '4.14'
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('D23').Value = '1.97'
$ws.Range('E23').Value = '  -6.21%  '
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').Value = '65.73'
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('D26').Value = '9.29'
$ws.Range('E26').Value = '  +3.44%  '
$ws.Range('D27').Value = '611.38'
$ws.Range('E27').Value = '  -1.55%  '
$ws.Range('D28').Value = '2.587.44'
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('D29').Value = '0.0₃0968'
$ws.Range('E29').Value = '  -2.87%  '
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('D31').Value = '1.44'
$ws.Range('E31').Value = '  -3.46%  '
$ws.Range('D32').Value = '8.02'
$ws.Range('E32').Value = '  -1.46%  '
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('E34').Value = '  +1.12%  '
$ws.Range('D35').Value = '4.94'
$ws.Range('E35').Value = '  -3.61%  '
$ws.Range('E36').Value = '  +0.28%  '
$ws.Range('D37').Value = '1.45'
$ws.Range('E37').Value = '  -2.96%  '
$ws.Range('D38').Value = '0.378'
$ws.Range('E38').Value = '  +0.38%  '
$ws.Range('D39').Value = '5.38'
$ws.Range('E39').Value = '  +1.05%  '
$ws.Range('D40').Value = '150.22'
$ws.Range('D41').Value = '18.48'
$ws.Range('E41').Value = '  -1.60%  '
$ws.Range('D42').Value = '1.74'
$ws.Range('E42').Value = '  -2.14%  '
$ws.Range('D43').Value = '42.83'
$ws.Range('E43').Value = '  +1.88%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').Value = '2.53'
$ws.Range('E45').Value = '  -2.13%  '
$ws.Range('D46').Value = '143.46'
$ws.Range('E46').Value = '  -2.06%  '
$ws.Range('D47').Value = '3.64'
$ws.Range('E47').Value = '  -3.24%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0252'
$ws.Range('E48').Value = '  +18.26%  '
$ws.Range('D49').Value = '0.0527'
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '0.607'
$ws.Range('E50').Value = '  +1.87%  '
$ws.Range('D51').Value = '19.77'
$ws.Range('E51').Value = '  -4.71%  '
